$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 51
$ws.Range("H51").Value = 7799.75
$ws.Range("I51").Value = 6774.5
$ws.Range("J51").Value = 8141.5
$ws.Range("K51").Value = 6774.5
$ws.Range("L51").Value = 8141.5
$ws.Range("M51").Value = -6290.5
$ws.Range("N51").Value = -9109.5

# Row 62
$ws.Range("H62").Value = 3900.037
$ws.Range("I62").Value = 2838.2104
$ws.Range("K62").Value = 2838.2104
$ws.Range("M62").Value = -2214.2104

# Row 65
$ws.Range("H65").Value = 3900.037
$ws.Range("I65").Value = 2838.2104
$ws.Range("K65").Value = 14191.052
$ws.Range("M65").Value = -11071.052

# Row 69
$ws.Range("H69").Value = 8071
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents() | Out-Null

# Row 72
$ws.Range("H72").Value = 8071
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents() | Out-Null

# Row 132
$ws.Range("H132").Value = 1060.5
$ws.Range("I132").Value = 996.5854
$ws.Range("K132").Value = 2989.7562
$ws.Range("M132").Value = -459.7562000000003

# Row 135
$ws.Range("H135").Value = 744.2308
$ws.Range("I135").Value = 606.6
$ws.Range("K135").Value = 5459.400000000001
$ws.Range("M135").Value = -2924.400000000001

# Row 138
$ws.Range("H138").Value = 2827.122
$ws.Range("J138").Value = 3221.2173
$ws.Range("L138").Value = 9663.651899999999
$ws.Range("N138").Value = -19943.6519

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4490.674
$ws.Range("I32").Value = 3698.6743
$ws.Range("J32").Value = 15842.667
$ws.Range("K32").Value = 3698.6743
$ws.Range("L32").Value = 15842.667
$ws.Range("M32").Value = -3411.6743
$ws.Range("N32").Value = -16416.667

# Row 74
$ws.Range("H74").Value = 37040296
$ws.Range("I74").Value = 55559444
$ws.Range("K74").Value = 55559444
$ws.Range("M74").Value = -55558570

# Row 77
$ws.Range("H77").Value = 37040296
$ws.Range("I77").Value = 55559444
$ws.Range("K77").Value = 277797220
$ws.Range("M77").Value = -277792852

# Row 124
$ws.Range("H124").Value = 42500
$ws.Range("J124").Value = 42500
$ws.Range("L124").Value = 42500
$ws.Range("N124").Value = -52320

# Row 132
$ws.Range("H132").Value = 2766.5557
$ws.Range("I132").Value = 2184.238
$ws.Range("K132").Value = 6552.714
$ws.Range("M132").Value = -4022.714

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 4295.2383
$ws.Range("I86").Value = 2661.4119
$ws.Range("J86").Value = 11239
$ws.Range("K86").Value = 2661.4119
$ws.Range("L86").Value = 11239
$ws.Range("M86").Value = -1538.4119
$ws.Range("N86").Value = -13485

# Row 89
$ws.Range("H89").Value = 4295.2383
$ws.Range("I89").Value = 2661.4119
$ws.Range("J89").Value = 11239
$ws.Range("K89").Value = 13307.0595
$ws.Range("L89").Value = 56195
$ws.Range("M89").Value = -7691.059499999999
$ws.Range("N89").Value = -67427

# Row 94
$ws.Range("H94").Value = 1277.8823
$ws.Range("I94").Value = 1261.6154
$ws.Range("J94").Value = 1330.75
$ws.Range("K94").Value = 1261.6154
$ws.Range("L94").Value = 1330.75
$ws.Range("M94").Value = -810.6153999999999
$ws.Range("N94").Value = -2232.75

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 570.5599999999999
$ws.Range("I7").Value = 742.8125
$ws.Range("K7").Value = 742.8125
$ws.Range("M7").Value = -629.8125

# Row 16
$ws.Range("H16").Value = 2065.9285
$ws.Range("I16").Value = 1744
$ws.Range("J16").Value = 3997.5
$ws.Range("K16").Value = 1744
$ws.Range("L16").Value = 3997.5
$ws.Range("M16").Value = -1457
$ws.Range("N16").Value = -4571.5

# Row 31
$ws.Range("H31").Value = 35850.668
$ws.Range("I31").Value = 4442.6665
$ws.Range("K31").Value = 4442.6665
$ws.Range("M31").Value = -4147.6665

# Row 34
$ws.Range("H34").Value = 35850.668
$ws.Range("I34").Value = 4442.6665
$ws.Range("K34").Value = 4442.6665
$ws.Range("M34").Value = -4240.6665

# Row 58
$ws.Range("H58").Value = 3622.52
$ws.Range("I58").Value = 1641.2354
$ws.Range("K58").Value = 1641.2354
$ws.Range("M58").Value = -1438.2354

# Row 62
$ws.Range("H62").Value = 6031.5
$ws.Range("I62").Value = 4228.8335
$ws.Range("K62").Value = 4228.8335
$ws.Range("M62").Value = -3604.8335

# Row 65
$ws.Range("H65").Value = 6031.5
$ws.Range("I65").Value = 4228.8335
$ws.Range("K65").Value = 21144.1675
$ws.Range("M65").Value = -18024.1675

# Row 113
$ws.Range("H113").Value = 2065.9285
$ws.Range("I113").Value = 1744
$ws.Range("J113").Value = 3997.5
$ws.Range("K113").Value = 1744
$ws.Range("L113").Value = 3997.5
$ws.Range("M113").Value = 426
$ws.Range("N113").Value = -8337.5

# Row 134
$ws.Range("H134").Value = 3825.6155
$ws.Range("I134").Value = 2201
$ws.Range("J134").Value = 7481
$ws.Range("K134").Value = 6603
$ws.Range("L134").Value = 22443
$ws.Range("M134").Value = -4068
$ws.Range("N134").Value = -27513

# Row 136
$ws.Range("H136").Value = 3622.52
$ws.Range("I136").Value = 1641.2354
$ws.Range("K136").Value = 4923.706200000001
$ws.Range("M136").Value = -2373.706200000001

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 54868.09
$ws.Range("J2").Value = 54868.09
$ws.Range("L2").Value = 329208.54
$ws.Range("N2").Value = -329434.54

# Row 5
$ws.Range("H5").Value = 2119.45
$ws.Range("I5").Value = 772.2727
$ws.Range("K5").Value = 2316.8181
$ws.Range("M5").Value = -2204.8181

# Row 131
$ws.Range("H131").Value = 11439641
$ws.Range("J131").Value = 24306710
$ws.Range("L131").Value = 72920130
$ws.Range("N131").Value = -72930210

# Row 132
$ws.Range("H132").Value = 3806.0588
$ws.Range("J132").Value = 4833.5
$ws.Range("L132").Value = 43501.5
$ws.Range("N132").Value = -48561.5

# Row 135
$ws.Range("H135").Value = 2119.45
$ws.Range("I135").Value = 772.2727
$ws.Range("K135").Value = 6950.454299999999
$ws.Range("M135").Value = -4415.454299999999

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 1221.2
$ws.Range("I2").Value = 31.5
$ws.Range("K2").Value = 31.5
$ws.Range("M2").Value = 81.5

# Row 32
$ws.Range("H32").Value = 44999.5
$ws.Range("J32").Value = 44999.5
$ws.Range("L32").Value = 44999.5
$ws.Range("N32").Value = -45591.5

# Row 70
$ws.Range("H70").Value = 14929.565
$ws.Range("I70").Value = 5552.0586
$ws.Range("K70").Value = 5552.0586
$ws.Range("M70").Value = -5282.0586

# Row 73
$ws.Range("H73").Value = 14929.565
$ws.Range("I73").Value = 5552.0586
$ws.Range("K73").Value = 5552.0586
$ws.Range("M73").Value = -4616.0586

# Row 122
$ws.Range("H122").Value = 3203.4
$ws.Range("I122").Value = 2597.5
$ws.Range("K122").Value = 7792.5
$ws.Range("M122").Value = -5342.5

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 3917.3635
$ws.Range("J46").Value = 5228.7144
$ws.Range("L46").Value = 5228.7144
$ws.Range("N46").Value = -5604.7144

# Row 68
$ws.Range("H68").Value = 4271.9546
$ws.Range("I68").Value = 2999
$ws.Range("K68").Value = 2999
$ws.Range("M68").Value = -2250

# Row 71
$ws.Range("H71").Value = 4271.9546
$ws.Range("I71").Value = 2999
$ws.Range("K71").Value = 14995
$ws.Range("M71").Value = -11251

# Row 122
$ws.Range("H122").Value = 8955.333000000001
$ws.Range("I122").Value = 7993.4287
$ws.Range("K122").Value = 23980.2861
$ws.Range("M122").Value = -21530.2861

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 3957.2727
$ws.Range("I81").Value = 2567.25
$ws.Range("J81").Value = 7664
$ws.Range("K81").Value = 5134.5
$ws.Range("L81").Value = 15328
$ws.Range("M81").Value = -4073.5
$ws.Range("N81").Value = -17450

# Row 84
$ws.Range("H84").Value = 3957.2727
$ws.Range("I84").Value = 2567.25
$ws.Range("J84").Value = 7664
$ws.Range("K84").Value = 25672.5
$ws.Range("L84").Value = 76640
$ws.Range("M84").Value = -20368.5
$ws.Range("N84").Value = -87248

# Row 113
$ws.Range("H113").Value = 1123.2
$ws.Range("J113").Value = 1123.2
$ws.Range("L113").Value = 3369.6
$ws.Range("N113").Value = -7709.6

# Row 122
$ws.Range("H122").Value = 7425.6553
$ws.Range("I122").Value = 2009.1875
$ws.Range("K122").Value = 6027.5625
$ws.Range("M122").Value = -3577.5625

# Row 138
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents() | Out-Null
